$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Placeholder"
$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "0"
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = "В работе"
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = "admin"
$ws.Cells.Item(5, 15).Value = "На месте"

# Row 6
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(6, 2).Value = "Placeholder"
$ws.Cells.Item(6, 3).NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = "2"
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = "В работе"
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = "admin"
$ws.Cells.Item(6, 15).Value = "Все просрано"
